$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Tabelle2 (chart source data): append the 2017 data row ---
[void]$ws2.Range("A7").Copy($ws2.Range("A8"))
$ws2.Range("A8").Value() = 2017
$ws2.Range("B8").Value() = 95
$ws2.Range("C8").Value() = 5
[void]$ws2.Range("C12").Select()

# --- Tabelle1: append matching 2017 row to the visible data table ---
[void]$ws1.Range("A11").Copy($ws1.Range("A12"))
$ws1.Range("A12").Value() = 2017
[void]$ws1.Range("B11").Copy($ws1.Range("B12"))
$ws1.Range("B12").Value() = "> 95"

# --- Tabelle1: update the WHO data-source note with the "as of" date ---
$ws1.Range("B35").Value() = "Weltgesundheitsorganisation (WHO), Stand 11.04.2019"

# --- Chart: extend the plotted ranges to include the new 2017 row ---
$co = $ws1.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()
$s1 = $sc.Item(1)
$s1.Formula = "=SERIES(,Tabelle2!`$A`$1:`$A`$8,Tabelle2!`$B`$1:`$B`$8,1)"
$s2 = $sc.Item(2)
$s2.Formula = "=SERIES(,Tabelle2!`$A`$1:`$A`$8,Tabelle2!`$C`$1:`$C`$8,2)"

# --- restore view/selection state on Tabelle1 ---
[void]$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
[void]$ws1.Range("J24").Select()
